# Insert a new row at 435, shifting existing rows 435:497 down to 436:498.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("435:435").Insert()

# Fill the newly inserted row 435 with the new data record.
$ws.Cells.Item(435, 1).Value = 7
$ws.Cells.Item(435, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(435, 3).Value = "Ñuble"
$ws.Cells.Item(435, 4).Value = 45142
$ws.Cells.Item(435, 5).Value = 16
$ws.Cells.Item(435, 6).Value = 100114013
$ws.Cells.Item(435, 7).Value = "Zanahoria"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Segunda"
$ws.Cells.Item(435, 10).Value = 150
$ws.Cells.Item(435, 11).Value = 5000
$ws.Cells.Item(435, 12).Value = 5000
$ws.Cells.Item(435, 13).Value = 5000
$ws.Cells.Item(435, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(435, 15).Value = "Región de Ñuble"
$ws.Cells.Item(435, 16).Value = 250
$ws.Cells.Item(435, 17).Value = 20
$ws.Cells.Item(435, 18).Value = "Hortaliza"
